$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "Tickers" header in Z1, matching the style of the other header cells (Y1)
$ws.Range("Y1").Copy() | Out-Null
$ws.Range("Z1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("Z1").Value = "Tickers"

# Append new experiment log rows (12-20 / Experiment IDs 11-19)
$ws.Range("A12").Value = 11
$ws.Range("B12").Value = 'Auto-log: Q=2, D=2, Skip=concat'
$ws.Range("C12").Value = 2
$ws.Range("D12").Value = 1
$ws.Range("E12").Value = 0
$ws.Range("F12").Value = 16
$ws.Range("G12").Value = 20
$ws.Range("H12").Value = $True
$ws.Range("I12").Value = 2
$ws.Range("J12").Value = 2
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 'concat'
$ws.Range("N12").Value = $False
$ws.Range("O12").Value = 0
$ws.Range("P12").Value = $False
$ws.Range("R12").Value = 15
$ws.Range("S12").Value = 2
$ws.Range("T12").Value = 0.0002755922186123725
$ws.Range("U12").Value = 0.008164557970303576
$ws.Range("V12").Value = 14.92195510864258
$ws.Range("W12").Value = 19.71026420593262
$ws.Range("X12").Value = 7.186328411102295
$ws.Range("Y12").Value = 'run with Quantum layer'
$ws.Range("A13").Value = 12
$ws.Range("B13").Value = 'Auto-log: Q=2, D=2, Skip=concat'
$ws.Range("C13").Value = 2
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 16
$ws.Range("G13").Value = 20
$ws.Range("H13").Value = $True
$ws.Range("I13").Value = 2
$ws.Range("J13").Value = 2
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 'concat'
$ws.Range("N13").Value = $False
$ws.Range("O13").Value = 0
$ws.Range("P13").Value = $False
$ws.Range("R13").Value = 15
$ws.Range("S13").Value = 2
$ws.Range("T13").Value = 0.0002311510873527917
$ws.Range("U13").Value = 0.006108876623329706
$ws.Range("V13").Value = 12.91889095306396
$ws.Range("W13").Value = 17.25859451293945
$ws.Range("X13").Value = 6.20987606048584
$ws.Range("Y13").Value = 'run with Quantum layer'
$ws.Range("A14").Value = 13
$ws.Range("B14").Value = 'Auto-log: Q=2, D=2, Skip=concat'
$ws.Range("C14").Value = 4
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 16
$ws.Range("G14").Value = 20
$ws.Range("H14").Value = $True
$ws.Range("I14").Value = 2
$ws.Range("J14").Value = 2
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 'concat'
$ws.Range("N14").Value = $False
$ws.Range("O14").Value = 0
$ws.Range("P14").Value = $False
$ws.Range("R14").Value = 15
$ws.Range("S14").Value = 2
$ws.Range("T14").Value = 0.0006707275329279669
$ws.Range("U14").Value = 0.01115978131565498
$ws.Range("V14").Value = 15.28087520599365
$ws.Range("W14").Value = 21.42845153808594
$ws.Range("X14").Value = 7.279820442199707
$ws.Range("Y14").Value = 'run with Quantum layer'
$ws.Range("A15").Value = 14
$ws.Range("B15").Value = 'Auto-log: Q=2, D=2, Skip=concat'
$ws.Range("C15").Value = 4
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 16
$ws.Range("G15").Value = 20
$ws.Range("H15").Value = $True
$ws.Range("I15").Value = 2
$ws.Range("J15").Value = 2
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 'concat'
$ws.Range("N15").Value = $False
$ws.Range("O15").Value = 0
$ws.Range("P15").Value = $False
$ws.Range("R15").Value = 15
$ws.Range("S15").Value = 2
$ws.Range("T15").Value = 0.00001132239412982017
$ws.Range("U15").Value = 0.00005341193536878563
$ws.Range("V15").Value = 0.1708307266235352
$ws.Range("W15").Value = 0.1843285113573074
$ws.Range("X15").Value = 0.60400390625
$ws.Range("Y15").Value = 'run with Quantum layer'
$ws.Range("A16").Value = 15
$ws.Range("B16").Value = 'Auto-log: Q=2, D=2, Skip=concat'
$ws.Range("C16").Value = 4
$ws.Range("D16").Value = 0
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 16
$ws.Range("G16").Value = 20
$ws.Range("H16").Value = $True
$ws.Range("I16").Value = 2
$ws.Range("J16").Value = 2
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 'concat'
$ws.Range("N16").Value = $False
$ws.Range("O16").Value = 0
$ws.Range("P16").Value = $False
$ws.Range("R16").Value = 15
$ws.Range("S16").Value = 2
$ws.Range("T16").Value = 0.0006707280517152939
$ws.Range("U16").Value = 0.01115977805784496
$ws.Range("V16").Value = 15.28087711334229
$ws.Range("W16").Value = 21.42845344543457
$ws.Range("X16").Value = 7.279821395874023
$ws.Range("Y16").Value = 'run with Quantum layer'
$ws.Range("A17").Value = 16
$ws.Range("B17").Value = 'Auto-log: Q=2, D=2, Skip=concat'
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 0
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 16
$ws.Range("G17").Value = 20
$ws.Range("H17").Value = $True
$ws.Range("I17").Value = 2
$ws.Range("J17").Value = 2
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 'concat'
$ws.Range("N17").Value = $False
$ws.Range("O17").Value = 0
$ws.Range("P17").Value = $False
$ws.Range("R17").Value = 15
$ws.Range("S17").Value = 2
$ws.Range("T17").Value = 0.0002620218473702095
$ws.Range("U17").Value = 0.0008594085127387711
$ws.Range("V17").Value = 3.520127534866333
$ws.Range("W17").Value = 4.569629192352295
$ws.Range("X17").Value = 2.406062841415405
$ws.Range("Y17").Value = 'run with Quantum layer'
$ws.Range("A18").Value = 17
$ws.Range("B18").Value = 'Auto-log: Q=2, D=2, Skip=concat'
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 0
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 16
$ws.Range("G18").Value = 20
$ws.Range("H18").Value = $True
$ws.Range("I18").Value = 2
$ws.Range("J18").Value = 2
$ws.Range("K18").Value = 3
$ws.Range("L18").Value = 'concat'
$ws.Range("N18").Value = $False
$ws.Range("O18").Value = 0
$ws.Range("P18").Value = $False
$ws.Range("R18").Value = 15
$ws.Range("S18").Value = 2
$ws.Range("T18").Value = 0.0002620218473702095
$ws.Range("U18").Value = 0.0008594085127387711
$ws.Range("V18").Value = 3.520127534866333
$ws.Range("W18").Value = 4.569629192352295
$ws.Range("X18").Value = 2.406062841415405
$ws.Range("Y18").Value = 'Multi-ticker experiment with LSTM+Quantum'
$ws.Range("Z18").Value = 'AAPL, MSFT, GOOGL'
$ws.Range("A19").Value = 18
$ws.Range("B19").Value = 'Auto-log: Q=2, D=2, Skip=concat'
$ws.Range("C19").Value = 4
$ws.Range("D19").Value = 0
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 16
$ws.Range("G19").Value = 20
$ws.Range("H19").Value = $True
$ws.Range("I19").Value = 2
$ws.Range("J19").Value = 2
$ws.Range("K19").Value = 3
$ws.Range("L19").Value = 'concat'
$ws.Range("N19").Value = $False
$ws.Range("O19").Value = 0
$ws.Range("P19").Value = $False
$ws.Range("R19").Value = 15
$ws.Range("S19").Value = 2
$ws.Range("T19").Value = 0.0002620218473702095
$ws.Range("U19").Value = 0.0008594085127387711
$ws.Range("V19").Value = 3.520127534866333
$ws.Range("W19").Value = 4.569629192352295
$ws.Range("X19").Value = 2.406062841415405
$ws.Range("Y19").Value = 'run with Quantum layer'
$ws.Range("Z19").Value = 'AAPL, MSFT, GOOGL'
$ws.Range("A20").Value = 19
$ws.Range("B20").Value = 'Auto-log: Q=2, D=2, Skip=concat'
$ws.Range("C20").Value = 4
$ws.Range("D20").Value = 0
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 16
$ws.Range("G20").Value = 20
$ws.Range("H20").Value = $True
$ws.Range("I20").Value = 2
$ws.Range("J20").Value = 2
$ws.Range("K20").Value = 3
$ws.Range("L20").Value = 'concat'
$ws.Range("N20").Value = $False
$ws.Range("O20").Value = 0
$ws.Range("P20").Value = $False
$ws.Range("R20").Value = 15
$ws.Range("S20").Value = 2
$ws.Range("T20").Value = 0.0002620216811455395
$ws.Range("U20").Value = 0.0008594098481092047
$ws.Range("V20").Value = 3.520127534866333
$ws.Range("W20").Value = 4.569628238677979
$ws.Range("X20").Value = 2.406063079833984
$ws.Range("Y20").Value = 'run with Quantum layer'
$ws.Range("Z20").Value = 'AAPL, MSFT, GOOGL'

Write-Host "Applied qml_experiment_log edits: added Tickers column and 9 new experiment rows."
